$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 7
$ws1.Range("H2").Value = 7.8

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 3
$ws2.Range("F2").Value = 20
$ws2.Range("G2").Value = 86.95999999999999
$ws2.Range("H2").Value = 7.8

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 3
$ws3.Range("F2").Value = 20
$ws3.Range("G2").Value = 86.95999999999999
$ws3.Range("H2").Value = 8.1

# --- Sheet "Rescatables" ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows.Item(3).Delete()
$ws4.Rows.Item(2).Delete()
